# Updates cryptos list price/volume(1h) values (and the three swapped
# coin rows 45-47) to match the refreshed snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Leading apostrophe forces Excel to store numeric-looking
    # strings (e.g. "1.00", "67.980.18") as text instead of
    # silently coercing them to numbers, matching the source
    # workbook where these columns are plain text cells.
    $range.Value = "'" + $value
    # Reset to the Normal style so the "number stored as text"
    # quote-prefix marker does not linger as an explicit style.
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.980.18"
Set-TextValue $ws.Range("E2") "  +0.94%  "
Set-TextValue $ws.Range("D3") "2.625.91"
Set-TextValue $ws.Range("E3") "  -0.16%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.14%  "
Set-TextValue $ws.Range("D5") "597.25"
Set-TextValue $ws.Range("E5") "  -0.77%  "
Set-TextValue $ws.Range("D6") "152.92"
Set-TextValue $ws.Range("E6") "  -0.05%  "
Set-TextValue $ws.Range("E7") "  -0.07%  "
Set-TextValue $ws.Range("D8") "0.545"
Set-TextValue $ws.Range("E8") "  -2.32%  "
Set-TextValue $ws.Range("D9") "2.624.71"
Set-TextValue $ws.Range("E9") "  -0.17%  "
Set-TextValue $ws.Range("E10") "  +8.22%  "
Set-TextValue $ws.Range("E11") "  -0.60%  "
Set-TextValue $ws.Range("D12") "5.20"
Set-TextValue $ws.Range("E12") "  +0.02%  "
Set-TextValue $ws.Range("E13") "  -1.45%  "
Set-TextValue $ws.Range("D14") "27.59"
Set-TextValue $ws.Range("E14") "  -0.53%  "
Set-TextValue $ws.Range("E15") "  +3.28%  "
Set-TextValue $ws.Range("D16") "3.098.15"
Set-TextValue $ws.Range("E16") "  -0.50%  "
Set-TextValue $ws.Range("D17") "67.839.11"
Set-TextValue $ws.Range("E17") "  +0.82%  "
Set-TextValue $ws.Range("D18") "2.615.90"
Set-TextValue $ws.Range("E18") "  -0.57%  "
Set-TextValue $ws.Range("D19") "372.34"
Set-TextValue $ws.Range("E19") "  +2.41%  "
Set-TextValue $ws.Range("D20") "11.27"
Set-TextValue $ws.Range("E20") "  +0.58%  "
Set-TextValue $ws.Range("D21") "7.47"
Set-TextValue $ws.Range("E21") "  -1.03%  "
Set-TextValue $ws.Range("E22") "  -1.56%  "
Set-TextValue $ws.Range("D23") "4.80"
Set-TextValue $ws.Range("E23") "  -2.49%  "
Set-TextValue $ws.Range("D24") "2.05"
Set-TextValue $ws.Range("E24") "  -3.22%  "
Set-TextValue $ws.Range("D25") "72.79"
Set-TextValue $ws.Range("E25") "  +9.59%  "
Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  -0.06%  "
Set-TextValue $ws.Range("D27") "9.87"
Set-TextValue $ws.Range("E27") "  -2.79%  "
Set-TextValue $ws.Range("E28") "  +1.81%  "
Set-TextValue $ws.Range("D30") "0.997"
Set-TextValue $ws.Range("E30") "  -0.38%  "
Set-TextValue $ws.Range("D31") "576.05"
Set-TextValue $ws.Range("E31") "  -0.34%  "
Set-TextValue $ws.Range("D32") "1.39"
Set-TextValue $ws.Range("E32") "  -0.30%  "
Set-TextValue $ws.Range("D33") "7.81"
Set-TextValue $ws.Range("E33") "  -0.41%  "
Set-TextValue $ws.Range("E34") "  -0.24%  "
Set-TextValue $ws.Range("E35") "  -0.07%  "
Set-TextValue $ws.Range("E36") "  -2.08%  "
Set-TextValue $ws.Range("D37") "1.51"
Set-TextValue $ws.Range("E37") "  -1.06%  "
Set-TextValue $ws.Range("D38") "159.01"
Set-TextValue $ws.Range("E38") "  +0.54%  "
Set-TextValue $ws.Range("D39") "19.16"
Set-TextValue $ws.Range("E39") "  -1.14%  "
Set-TextValue $ws.Range("E40") "  +4.16%  "
Set-TextValue $ws.Range("E41") "  -0.33%  "
Set-TextValue $ws.Range("D42") "5.30"
Set-TextValue $ws.Range("E42") "  +0.41%  "
Set-TextValue $ws.Range("E43") "  +1.24%  "
Set-TextValue $ws.Range("D44") "17.09"
Set-TextValue $ws.Range("E44") "  +4.53%  "
Set-TextValue $ws.Range("B45") "BabyDogeCoin"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D45") "0.0₆0312"
Set-TextValue $ws.Range("E45") "  +8.10%  "
Set-TextValue $ws.Range("B46") "USDe"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D46") "0.999"
Set-TextValue $ws.Range("E46") "  -0.03%  "
Set-TextValue $ws.Range("B47") "OKB"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D47") "40.41"
Set-TextValue $ws.Range("E47") "  -1.95%  "
Set-TextValue $ws.Range("D48") "155.10"
Set-TextValue $ws.Range("E48") "  -0.46%  "
Set-TextValue $ws.Range("D49") "3.69"
Set-TextValue $ws.Range("E49") "  -1.02%  "
Set-TextValue $ws.Range("E50") "  -1.89%  "
Set-TextValue $ws.Range("D51") "0.0780"
Set-TextValue $ws.Range("E51") "  -1.59%  "
